# Commit 7 For Rewriiten Code(made small changes 2)
# Update the "CorrectPoints" column (G) for rows 2-6 from 5 to 1,
# adjust the auto-fit row heights for rows 4 and 5 that shift as a
# result, and move the viewport/selection to B1 (topLeft) / G7 (active cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CorrectPoints column (G2:G6): 5 -> 1 ---
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1

# --- Row height adjustments (rows re-wrapped after edit) ---
$ws.Rows.Item(4).RowHeight = 82.05
$ws.Rows.Item(5).RowHeight = 95.5

# --- View state: scroll window so column B is left-most visible, select G7 ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("G7").Select()
